$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Bonjour1): extend weekly data across Feb columns (G:N) and tweak E4/F4
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = 10
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 10
$ws.Range("K4").Value = 10
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1
$ws.Range("N4").Value = 1

# Row 5 (Bonjour2): same value (2) across the new columns
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 2
$ws.Range("M5").Value = 2
$ws.Range("N5").Value = 2

# Row 6 (Muchacho1): same value (3) across the new columns
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 3
$ws.Range("I6").Value = 3
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 3
$ws.Range("M6").Value = 3
$ws.Range("N6").Value = 3

# Row 7 (Muchacho2): same value (4) across the new columns
$ws.Range("F7").Value = 4
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = 4
$ws.Range("I7").Value = 4
$ws.Range("J7").Value = 4
$ws.Range("K7").Value = 4
$ws.Range("L7").Value = 4
$ws.Range("M7").Value = 4
$ws.Range("N7").Value = 4

# Move the selection/active cell to match the author's final cursor position
$ws.Range("O4:P7").Select()
